$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "ProfileName"
$ws.Range("C2").Value = " admin KL "

$ws.Range("C4").Select()
